# Update the "保險" (insurance) worksheet (sheet4) with the full set of
# columns/rows describing each insurance policy record, matching the
# expanded schema already used on the other sheets (company/name/owner/
# category/date/legislator_name/legislator_id/source_file/index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# Start from a clean slate for the region we are rewriting.
$ws.Range("A1:K9").ClearContents()

# Header row
$ws.Cells.Item(1,2).Value = "company"
$ws.Cells.Item(1,3).Value = "name"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "property_category"
$ws.Cells.Item(1,6).Value = "category"
$ws.Cells.Item(1,7).Value = "date"
$ws.Cells.Item(1,8).Value = "legislator_name"
$ws.Cells.Item(1,9).Value = "legislator_id"
$ws.Cells.Item(1,10).Value = "source_file"
$ws.Cells.Item(1,11).Value = "index"

# Data rows: A=index, B=company, C=name, D=owner, E=property_category,
# F=category, G=date, H=legislator_name, I=legislator_id, J=source_file, K=index
$rows = @(
  @(93,  "南山人壽", "子女教育保險",             "廖述嘉", "insurance", "normal", "2012-02-10", "盧秀燕", 869, "tmp61a71", 93),
  @(94,  "中華郵政", "十年快樂兒童還本終身壽險", "盧秀燕", "insurance", "normal", "2012-02-10", "盧秀燕", 869, "tmp61a71", 94),
  @(95,  "中華郵政", "十年快樂兒童還本終身壽險", "盧秀燕", "insurance", "normal", "2012-02-10", "盧秀燕", 869, "tmp61a71", 95),
  @(96,  "中華郵政", "金寶貝兒童保險",           "盧秀燕", "insurance", "normal", "2012-02-10", "盧秀燕", 869, "tmp61a71", 96),
  @(97,  "中國人壽", "得意人生終身保險",         "盧秀燕", "insurance", "normal", "2012-02-10", "盧秀燕", 869, "tmp61a71", 97),
  @(98,  "中國人壽", "得意人生終身保險",         "盧秀燕", "insurance", "normal", "2012-02-10", "盧秀燕", 869, "tmp61a71", 98),
  @(99,  "中國人壽", "得意人生終身保險",         "盧秀燕", "insurance", "normal", "2012-02-10", "盧秀燕", 869, "tmp61a71", 99),
  @(100, "中國人壽", "得意人生終身保險",         "盧秀燕", "insurance", "normal", "2012-02-10", "盧秀燕", 869, "tmp61a71", 100)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value  = $row[0]
    $ws.Cells.Item($r,2).Value  = $row[1]
    $ws.Cells.Item($r,3).Value  = $row[2]
    $ws.Cells.Item($r,4).Value  = $row[3]
    $ws.Cells.Item($r,5).Value  = $row[4]
    $ws.Cells.Item($r,6).Value  = $row[5]
    $ws.Cells.Item($r,7).Value  = $row[6]
    $ws.Cells.Item($r,8).Value  = $row[7]
    $ws.Cells.Item($r,9).Value  = $row[8]
    $ws.Cells.Item($r,10).Value = $row[9]
    $ws.Cells.Item($r,11).Value = $row[10]
    $r = $r + 1
}

Write-Output "sheet4 (insurance) updated: $($r-2) data rows written"
